$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: N2 gets style 17 (same as L2/M2), empty value ---
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

# --- Row 3: N3 gets style 13 (same as L3/M3), value 2022 ---
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2022

# --- Row 4: N4 gets style 14 (same as M4), value 11.927942610539198 ---
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 11.927942610539198

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 3.0909744679837434
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 20.963679772397647
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 4.6002717699014832
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 0
$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 9.112830865859129
$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 3.5391993253978327
$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 0.30955295909412422
$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 6.73157537222552
$ws.Range("M13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 3.9173330796393815
$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 0.7444796831494469
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = 7.104530072727953
$ws.Range("M16").Copy()
$ws.Range("N16").PasteSpecial(-4122)
$ws.Range("N16").Value = 23.0957399744971
$ws.Range("M17").Copy()
$ws.Range("N17").PasteSpecial(-4122)
$ws.Range("N17").Value = 2.6274648905004008
$ws.Range("M18").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("N18").Value = 43.176223433734158
$ws.Range("M19").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("N19").Value = 7.6660105666632132
$ws.Range("M20").Copy()
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("N20").Value = 0.83437630371297455
$ws.Range("M21").Copy()
$ws.Range("N21").PasteSpecial(-4122)
$ws.Range("N21").Value = 14.406256431364477
$ws.Range("M22").Copy()
$ws.Range("N22").PasteSpecial(-4122)
$ws.Range("N22").Value = 34.201612992199827
$ws.Range("M23").Copy()
$ws.Range("N23").PasteSpecial(-4122)
$ws.Range("N23").Value = 4.4521615244201058
$ws.Range("M24").Copy()
$ws.Range("N24").PasteSpecial(-4122)
$ws.Range("N24").Value = 63.433733622066185
$ws.Range("M25").Copy()
$ws.Range("N25").PasteSpecial(-4122)
$ws.Range("N25").Value = 20.535408979625672
$ws.Range("M26").Copy()
$ws.Range("N26").PasteSpecial(-4122)
$ws.Range("N26").Value = 7.8632542639432348
$ws.Range("M27").Copy()
$ws.Range("N27").PasteSpecial(-4122)
$ws.Range("N27").Value = 33.368028499329796
$ws.Range("M28").Copy()
$ws.Range("N28").PasteSpecial(-4122)
$ws.Range("N28").Value = 19.301652062045072
$ws.Range("M29").Copy()
$ws.Range("N29").PasteSpecial(-4122)
$ws.Range("N29").Value = 7.1220113855063829
$ws.Range("M30").Copy()
$ws.Range("N30").PasteSpecial(-4122)
$ws.Range("N30").Value = 34.008685896558866
$ws.Range("M31").Copy()
$ws.Range("N31").PasteSpecial(-4122)
$ws.Range("N31").Value = 7.8668258762379715
$ws.Range("M32").Copy()
$ws.Range("N32").PasteSpecial(-4122)
$ws.Range("N32").Value = 1.7266187050359711

# --- Row 33: N33 gets style 16 (same as M33), value 13.723068478111704 ---
$ws.Range("M33").Copy()
$ws.Range("N33").PasteSpecial(-4122)
$ws.Range("N33").Value = 13.723068478111704

# --- Row 34: N34 gets a brand-new style (fontId=1 Times New Roman 11, no border,
#     no alignment, no number format). Build it by starting from a no-alignment
#     donor cell (L3, style 13: bold 9pt Times New Roman + bottom border), then
#     adjust the font to regular 11pt and clear the border. Each of these Font
#     sub-property edits (Bold / Size) preserves "no alignment" as long as the
#     starting xf had no alignment to begin with.
$ws.Range("L3").Copy()
$ws.Range("N34").PasteSpecial(-4122)
$ws.Range("N34").Font.Bold = $false
$ws.Range("N34").Font.Size = 11
$ws.Range("N34").Borders.LineStyle = -4142

# --- Update the active selection to match the target workbook state ---
$ws.Range("O6").Select()
